$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 5) logging another meal entry.
# Column layout: A=DateTime (Eat start), B=Home food?, C=Foods,
# D=Did not eat myself, E=Cost, F=My rating (/5), G=Comments
$ws.Range("A5").Value = 44602.463287037041
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "idli-sambhar-chutney"
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = "After being a regular breakfast once/twice a week the speciality goes down."

# Comments column wraps text; match row 4's taller row height for the new row.
$ws.Rows("5").RowHeight = 28.8

# Leave the selection where the author ended up after entering the data.
$ws.Range("G6").Select() | Out-Null
